$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header: copy format from E1 (bold/bordered header style) then set text
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 6).Value = "time_taken"

# F2:F52 time_taken values (plain text, default/unstyled like surrounding data cells)
$timestamps = @(
    "2021-10-05 10:50:58.613003",
    "2021-10-05 10:50:58.613015",
    "2021-10-05 10:50:58.613019",
    "2021-10-05 10:50:58.613022",
    "2021-10-05 10:50:58.613026",
    "2021-10-05 10:50:58.613029",
    "2021-10-05 10:50:58.613032",
    "2021-10-05 10:50:58.613035",
    "2021-10-05 10:50:58.613038",
    "2021-10-05 10:50:58.613041",
    "2021-10-05 10:50:58.613044",
    "2021-10-05 10:50:58.613047",
    "2021-10-05 10:50:58.613050",
    "2021-10-05 10:50:58.613053",
    "2021-10-05 10:50:58.613056",
    "2021-10-05 10:50:58.613059",
    "2021-10-05 10:50:58.613062",
    "2021-10-05 10:50:58.613066",
    "2021-10-05 10:50:58.613069",
    "2021-10-05 10:50:58.613072",
    "2021-10-05 10:50:58.613075",
    "2021-10-05 10:50:58.613078",
    "2021-10-05 10:50:58.613081",
    "2021-10-05 10:50:58.613084",
    "2021-10-05 10:50:58.613087",
    "2021-10-05 10:50:58.613090",
    "2021-10-05 10:50:58.613093",
    "2021-10-05 10:50:58.613096",
    "2021-10-05 10:50:58.613099",
    "2021-10-05 10:50:58.613102",
    "2021-10-05 10:50:58.613105",
    "2021-10-05 10:50:58.613108",
    "2021-10-05 10:50:58.613112",
    "2021-10-05 10:50:58.613115",
    "2021-10-05 10:50:58.613118",
    "2021-10-05 10:50:58.613121",
    "2021-10-05 10:50:58.613124",
    "2021-10-05 10:50:58.613127",
    "2021-10-05 10:50:58.613130",
    "2021-10-05 10:50:58.613132",
    "2021-10-05 10:50:58.613136",
    "2021-10-05 10:50:58.613139",
    "2021-10-05 10:50:58.613142",
    "2021-10-05 10:50:58.613145",
    "2021-10-05 10:50:58.613148",
    "2021-10-05 10:50:58.613151",
    "2021-10-05 10:50:58.613154",
    "2021-10-05 10:50:58.613157",
    "2021-10-05 10:50:58.613160",
    "2021-10-05 10:50:58.613163",
    "2021-10-05 10:50:58.613166"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

